$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace the placeholder sample data with the first real attendee.
$ws.Range("B2").Value = "Tan Zhang En"
$ws.Range("C2").Value = "zhangen69@gmail.com"
$ws.Range("D2").Value = "Male"
$ws.Range("E2").Value = "940830-01-6651"
$ws.Range("F2").Value = "019-7765290"

# Columns G-J of row 2 (Organization/Address/Group/Remarks samples) are no
# longer used for this attendee.
$ws.Range("G2:J2").ClearContents()

# Row 3: second attendee.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Soo De Xiang"
$ws.Range("C3").Value = "dexiang@gmail.com"
$ws.Range("D3").Value = "Male"
$ws.Range("E3").Value = "980122-01-6412"
$ws.Range("F3").Value = "012-1325418"

# Turn the two email cells into mailto hyperlinks (adds the "Hyperlink" cell
# style automatically).
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:zhangen69@gmail.com", "", "", "zhangen69@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:dexiang@gmail.com", "", "", "dexiang@gmail.com")

# Move the active selection.
$ws.Range("B4").Select() | Out-Null
